$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YES PHN")

$ws.Range("AH2").Value = 99
$ws.Range("AI2").Value = 99
$ws.Range("AH3").Value = 99
$ws.Range("AI3").Value = 99
